$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 520
$ws.Range("B9").Value = 4

$ws.Range("C12").Select()
